$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 for the E99 communication-error entry
$ws.Rows(8).Insert()

$ws.Range("C8").Value = "E99"
$ws.Range("D8").Value = '통신 이상'
$ws.Range("E8").Value = '통신 케이블 확인'
$ws.Range("C8:E8").NumberFormat = "@"

# Re-apply the text number format to the numeric-code cells in column C
# that shifted down one row because of the insert above
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C19:C23").NumberFormat = "@"

# Append two new rows describing aligner setup / execution
$ws.Range("C24").Value = '얼라인 설정'
$ws.Range("C24").NumberFormat = "@"
$ws.Range("D24").Value = '1) 웨이퍼 크기 설정
[ WFS 6<CR><LF> ] (6”(100mm) 웨이퍼로 설정)
[ WFS 200<CR><LF> ] (200mm(8”) 웨이퍼로 설정)
2)  웨이퍼 타입 설정
[ WFT 1<CR><LF> ] (노치타입 웨이퍼로 설정)
[ WFT 2<CR><LF> ] (플랫타입 웨이퍼로 설정)
'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").WrapText = $true
$ws.Rows(24).RowHeight = 121.8

$ws.Range("C25").Value = '얼라인 실행'
$ws.Range("C25").NumberFormat = "@"
$ws.Range("D25").Value = '1) 얼라이너에 24V 전원 인가.
2)  [ORG] 명령 전송
주의 : 전원 인가 후 ORG 명령은 반드시 실행되어야 합니다.
      ORG명령은 웨이퍼가 올려져 있지 않은 상태에서 실행되어야 합니다.
3)  [DWL] 명령으로 변수 설정
형식 : DWL #,XXX,YYY,RRR
      # : 변수 번호 (1~10)
      XXX : 얼라이닝 과정 후의 센터링 위치의 X축 OFFSET 값 (단위 : mm)
   YYY : 얼라이닝 과정 후의 센터링 위치의 Y축 OFFSET 값 (단위 : mm)
      RRR : 얼라이닝 과정 후 회전축의 회전 양 (단위 : 0.1°)
4)  WTR이 웨이퍼를 PUT
주의 : 얼라이너에는 센터링 위치 결정 기능이 있지만, CCD센서가 최적의 성능을 발휘하기 위해, 웨이퍼의 중심과     
얼라이너 회전축의 오차는 5mm이내가 되어야 합니다. 
웨이퍼의 중심이 회전축과 5mm이상 멀어지게 되면 CCD센서가 정확한 값을 읽지 못할 수도 있으므로 WTR이
최대한 얼라이너의 중심에 웨이퍼를 올려놓도록 티칭해 주시기 바랍니다.      
5)  [ALG] 명령(얼라이닝 작업만 수행) 또는 [ALS #] 명령(얼라이닝 작업 후 ‘#’번호의 변수에 저장된 위치만큼 이동&회전)
6)  WTR이 웨이퍼를 GET.
7)  [RST] 명령을 송신해서 작업 대기 상태로 복귀. 4번 스텝부터 계속 작업 진행.
'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").WrapText = $true
$ws.Rows(25).RowHeight = 409.6

# Restore the original selection/view anchor
$ws.Range("D25").Select()

Write-Host "edit complete"
